$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 3841 (entire row), matching the selection recorded in the diff,
# then delete it so all the rows below shift up by one.
$ws.Rows.Item(3841).Select()
$ws.Rows.Item(3841).Delete()

# Restore the view/selection state recorded after the edit.
$ws.Rows.Item(3841).Select()
$excel.ActiveWindow.ScrollRow = 3818
